$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(11703,10335,9610,9610,9610,9610,9610,9610,9610,9610,9610,9610,9610,8989,8989,8989,8989,8989,8989,8989,8989,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,8724,7657,7657,7657,7657,7657,7657,7657,7657,7657,7657,7657,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573,7573)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
